$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3 corrections (Appenzeller-Herzog 2019 - van Dis 2020 relevance markers)
$ws.Range("E3").Value = 1
$ws.Range("H3").Value = 0.6591805433829974
$ws.Range("I3").Value = 0.07158544365920691
$ws.Range("K3").Value = 684.7777777777778
$ws.Range("Q3").Value = 66
$ws.Range("R3").Value = 86
$ws.Range("S3").Value = 256
$ws.Range("T3").Value = 715
$ws.Range("U3").Value = 1639
$ws.Range("V3").Value = 8990
$ws.Range("W3").Value = 8970
$ws.Range("X3").Value = 8800
$ws.Range("Y3").Value = 8341
$ws.Range("Z3").Value = 7417
$ws.Range("AF3").Value = 0.992712
$ws.Range("AG3").Value = 0.9905040000000001
$ws.Range("AH3").Value = 0.971731
$ws.Range("AI3").Value = 0.9210469999999999
$ws.Range("AJ3").Value = 0.819015

$wb.Save()
